$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop all existing hyperlinks up front; rows 6-14 (and their hyperlinks) are
# being removed entirely, and rows 2-5 get brand-new target URLs, so the
# whole hyperlink collection gets rebuilt from scratch below.
$ws.Range("A1").Hyperlinks.Delete()

# New row 2 <- data that used to live in row 4 (re-timestamped)
$ws.Range("A2").Value = "2025-10-12 06:23:30"
$ws.Range("B2").Value = "急募 PR Zoom/Meet×TLDV×ChatGPT×Notion×Slack 議事録ワークフロー構築依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5410688"
$ws.Range("G2").Value = 323
$ws.Range("H2").Value = "🔥GPT,ChatGPT"

# New row 3 <- data that used to live in row 7 (re-timestamped)
$ws.Range("A3").Value = "2025-10-12 06:23:30"
$ws.Range("B3").Value = "【自動売買】Excelと楽天RSSを活用したシステム開発依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5411684"
$ws.Range("G3").Value = 110
$ws.Range("H3").Value = "◆開発,システム開発"

# New row 4 <- data that used to live in row 8 (re-timestamped)
$ws.Range("A4").Value = "2025-10-12 06:23:30"
$ws.Range("B4").Value = "Laravelでのバックエンド開発:管理画面機能やDB管理・ポイント機能などの開発【フルリモート】"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5411736"
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = "◆開発 ◇管理"

# New row 5 <- data that used to live in row 12 (re-timestamped)
$ws.Range("A5").Value = "2025-10-12 06:23:30"
$ws.Range("B5").Value = "【急募】教育系のWEBサイトの作成"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5411679"
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = "◇サイト"

# Rows 6-14 no longer exist in the refreshed export - drop them, which also
# shrinks the sheet dimension down to A1:H5.
$ws.Rows("6:14").Delete()

# Column D/H got a bit narrower in this export.
$ws.Columns("D").ColumnWidth = 27.166666666666668
$ws.Columns("H").ColumnWidth = 13.166666666666666

# Re-create the hyperlinks for the four surviving data rows against their
# (possibly new) URLs.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5410688")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5411684")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5411736")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5411679")
